$d = $word.ActiveDocument

# Replace the FILLER placeholder with the first paragraph of real content.
$d.Content.Find.Execute("FILLER", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "The Florida Polytechnic University SGA Department of NOVA is comprised of the SGA Director of NOVA and additional deputies, such as the Deputy of Logistics, Deputy if Marketing, and the Deputy of Finance.",
                         2)

# Append the remaining four paragraphs after the first one.
$p1 = $d.Paragraphs(1).Range
$p1.InsertParagraphAfter()

$p2 = $d.Paragraphs(2).Range
$p2.InsertAfter("The Department of NOVA is responsible for the planning and overseeing of the annual NOVA")
$p2.InsertParagraphAfter()

$p3 = $d.Paragraphs(3).Range
$p3.InsertAfter("banquet meant to celebrate the students and staff of Florida Polytechnic University and their")
$p3.InsertParagraphAfter()

$p4 = $d.Paragraphs(4).Range
$p4.InsertAfter("accomplishments.")
$p4.InsertParagraphAfter()

$p5 = $d.Paragraphs(5).Range
$p5.InsertAfter("Contact: SGA- nova@floridapoly.edu")
